$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.797565406020876
$ws.Range("C2").Value = 3.954677509319504
$ws.Range("B3").Value = 4.677536928975535
$ws.Range("C3").Value = 8.848840784703315
$ws.Range("B4").Value = 5.961097891796423
$ws.Range("C4").Value = 13.1206480759472
$ws.Range("B5").Value = 8.440282872316885
$ws.Range("C5").Value = 17.02725871023842
$ws.Range("B6").Value = 10.95498289802889
$ws.Range("C6").Value = 20.83818187421087
$ws.Range("B7").Value = 13.07875922300305
$ws.Range("C7").Value = 25.29738175455026
$ws.Range("B8").Value = 17.80016034760816
$ws.Range("C8").Value = 29.56586737311416
$ws.Range("B9").Value = 19.28567799416857
$ws.Range("C9").Value = 33.72965533858511
$ws.Range("B10").Value = 21.387958567855
$ws.Range("C10").Value = 37.55415423775222
$ws.Range("B11").Value = 23.40313956480723
$ws.Range("C11").Value = 42.14128527319716
$ws.Range("B12").Value = 27.75961736467185
$ws.Range("C12").Value = 46.03924457173056
$ws.Range("B13").Value = 29.56011895275808
$ws.Range("C13").Value = 50.61062440751094
$ws.Range("B14").Value = 32.39816714436071
$ws.Range("C14").Value = 55.12703792155841
$ws.Range("B15").Value = 34.2334642766857
$ws.Range("C15").Value = 58.95976806748047
$ws.Range("B16").Value = 36.81098300271012
$ws.Range("C16").Value = 63.25057093536375
$ws.Range("B17").Value = 38.0869738536263
$ws.Range("C17").Value = 67.20549717436434
$ws.Range("B18").Value = 39.01040223842817
$ws.Range("C18").Value = 71.61082821570498
$ws.Range("B19").Value = 40.58517080948999
$ws.Range("C19").Value = 76.04381967141978
$ws.Range("B20").Value = 40.96499395184869
$ws.Range("C20").Value = 79.97588142188242
$ws.Range("B21").Value = 42.49981855189498
$ws.Range("C21").Value = 84.45035574612109
$ws.Range("B22").Value = 44.24553300055891
$ws.Range("C22").Value = 89.08401748825622
$ws.Range("B23").Value = 46.56138860543905
$ws.Range("C23").Value = 94.49465243399112
$ws.Range("B24").Value = 47.57029354010074
$ws.Range("C24").Value = 98.22162261075368
$ws.Range("B25").Value = 49.17519514688801
$ws.Range("C25").Value = 102.7885162915794
$ws.Range("B26").Value = 52.52708031835691
$ws.Range("C26").Value = 106.6912196580947
$ws.Range("B27").Value = 55.3121697002949
$ws.Range("C27").Value = 110.6380128416091
$ws.Range("B28").Value = 56.03913088202432
$ws.Range("C28").Value = 114.807946201576
$ws.Range("B29").Value = 58.59779263178197
$ws.Range("C29").Value = 119.0832017526606
$ws.Range("B30").Value = 59.2802852016593
$ws.Range("C30").Value = 123.051809335653
$ws.Range("B31").Value = 61.15484052803086
$ws.Range("C31").Value = 126.9124445608208
$ws.Range("B32").Value = 62.55796789195843
$ws.Range("C32").Value = 131.1917152746573
$ws.Range("B33").Value = 65.50024433760166
$ws.Range("C33").Value = 135.8347278042092
$ws.Range("B34").Value = 66.43556713155643
$ws.Range("C34").Value = 139.8388339900565
$ws.Range("B35").Value = 71.14613256641478
$ws.Range("C35").Value = 143.9474753276745
$ws.Range("B36").Value = 72.89303389590836
$ws.Range("C36").Value = 148.0625914835387
$ws.Range("B37").Value = 76.02930088495714
$ws.Range("C37").Value = 152.5602874393692
$ws.Range("B38").Value = 76.94635775713127
$ws.Range("C38").Value = 156.61346376025
$ws.Range("B39").Value = 78.04401230648861
$ws.Range("C39").Value = 161.2538625235578
$ws.Range("B40").Value = 79.98556704833244
$ws.Range("C40").Value = 165.6711942496314
$ws.Range("B41").Value = 81.80685047373714
$ws.Range("C41").Value = 170.0840128240419
$ws.Range("B42").Value = 83.47859264513181
$ws.Range("C42").Value = 174.2545901457182
$ws.Range("B43").Value = 88.4958612499025
$ws.Range("C43").Value = 178.6038153202076
$ws.Range("B44").Value = 90.81391240434517
$ws.Range("C44").Value = 182.6069770014621
$ws.Range("B45").Value = 94.58921014283855
$ws.Range("C45").Value = 187.07179896263
$ws.Range("B46").Value = 95.56702402625771
$ws.Range("C46").Value = 191.2239541815437
$ws.Range("B47").Value = 96.39900213424353
$ws.Range("C47").Value = 194.9291902859382
$ws.Range("B48").Value = 97.4779173062874
$ws.Range("C48").Value = 199.7390641106491
$ws.Range("B49").Value = 98.50031392032066
$ws.Range("C49").Value = 203.8977958260083
